$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.254.84'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.651.16'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0629'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.26'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.89%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.883.80'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.636.59'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.539'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.01'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.88%  '
$ws.Range("D17").Value = '27.241.18'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.02'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.49'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.27'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.20'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.40'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.85'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0508'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("D35").Value = '1.272.33'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("E37").Value = '  +2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.542'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.845'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.811'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.40'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").Value = '1.794.07'
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.16'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.27'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").Value = '  +15.73%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.73'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.55%  '
$ws.Range("E51").Value = '  +0.74%  '
